# Auto-generated Excel COM-interop script
# Update automàtic: dades i banners [2026-02-10 21:50]
#
# Refreshes the per-row DATA_EXTRACCIO scrape timestamp (col E) and the
# small numeric drift picked up on this run for a handful of stations
# (humidity %, pressure hPa, rainfall mm, wind gust, radiation, temps).
#
# Percent-looking values ("79%", "90%", ...) are written with a leading
# apostrophe so Excel keeps them as literal text (matching the source
# file, which stores every reading as a plain string) instead of auto-
# converting them into a numeric percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("E2").Value = '2026-02-10 21:48:18'

# Row 3
$ws.Range("E3").Value = '2026-02-10 21:48:20'

# Row 4
$ws.Range("E4").Value = '2026-02-10 21:48:22'
$ws.Range("H4").Value = "'79%"
$ws.Range("J4").Value = '1003.7 hPa'

# Row 5
$ws.Range("E5").Value = '2026-02-10 21:48:25'

# Row 6
$ws.Range("E6").Value = '2026-02-10 21:48:27'

# Row 7
$ws.Range("E7").Value = '2026-02-10 21:48:29'
$ws.Range("J7").Value = '1004.3 hPa'
$ws.Range("O7").Value = '15.4 °C'

# Row 8
$ws.Range("E8").Value = '2026-02-10 21:48:32'
$ws.Range("J8").Value = '1004.2 hPa'
$ws.Range("O8").Value = '12.3 °C'

# Row 9
$ws.Range("E9").Value = '2026-02-10 21:48:34'
$ws.Range("I9").Value = '3.4 mm'
$ws.Range("L9").Value = '15.5 km/h - 239º 21:13 TU'
$ws.Range("O9").Value = '9.0 °C'

# Row 10
$ws.Range("E10").Value = '2026-02-10 21:48:37'
$ws.Range("L10").Value = '14.4 km/h - 32º 21:14 TU'
$ws.Range("O10").Value = '10.4 °C'

# Row 11
$ws.Range("E11").Value = '2026-02-10 21:48:39'
$ws.Range("H11").Value = "'90%"

# Row 12
$ws.Range("E12").Value = '2026-02-10 21:48:41'
$ws.Range("H12").Value = "'98%"
$ws.Range("I12").Value = '4.3 mm'
$ws.Range("O12").Value = '9.2 °C'

# Row 13
$ws.Range("E13").Value = '2026-02-10 21:48:43'
$ws.Range("J13").Value = '1006.5 hPa'

# Row 14
$ws.Range("E14").Value = '2026-02-10 21:48:46'

# Row 15
$ws.Range("E15").Value = '2026-02-10 21:48:48'

# Row 16
$ws.Range("E16").Value = '2026-02-10 21:48:51'
$ws.Range("I16").Value = '26.6 mm'

# Row 17
$ws.Range("E17").Value = '2026-02-10 21:48:53'
$ws.Range("H17").Value = "'94%"
$ws.Range("I17").Value = '0.3 mm'

# Row 18
$ws.Range("E18").Value = '2026-02-10 21:48:55'

# Row 19
$ws.Range("E19").Value = '2026-02-10 21:48:58'
$ws.Range("O19").Value = '6.9 °C'

# Row 20
$ws.Range("E20").Value = '2026-02-10 21:49:00'
$ws.Range("I20").Value = '11.6 mm'
$ws.Range("O20").Value = '0.6 °C'

# Row 21
$ws.Range("E21").Value = '2026-02-10 21:49:03'

# Row 22
$ws.Range("E22").Value = '2026-02-10 21:49:05'

# Row 23
$ws.Range("E23").Value = '2026-02-10 21:49:08'

# Row 24
$ws.Range("E24").Value = '2026-02-10 21:49:10'
$ws.Range("J24").Value = '1005.8 hPa'

# Row 25
$ws.Range("E25").Value = '2026-02-10 21:49:12'

# Row 26
$ws.Range("E26").Value = '2026-02-10 21:49:15'

# Row 27
$ws.Range("E27").Value = '2026-02-10 21:49:17'

# Row 28
$ws.Range("E28").Value = '2026-02-10 21:49:20'

# Row 29
$ws.Range("E29").Value = '2026-02-10 21:49:22'
$ws.Range("I29").Value = '0.7 mm'
$ws.Range("O29").Value = '10.8 °C'

# Row 30
$ws.Range("E30").Value = '2026-02-10 21:49:25'
$ws.Range("J30").Value = '1004.1 hPa'
$ws.Range("O30").Value = '9.4 °C'

# Row 31
$ws.Range("E31").Value = '2026-02-10 21:49:27'

# Row 32
$ws.Range("E32").Value = '2026-02-10 21:49:30'
$ws.Range("O32").Value = '10.7 °C'

# Row 33
$ws.Range("E33").Value = '2026-02-10 21:49:32'
$ws.Range("J33").Value = '1006.3 hPa'
$ws.Range("O33").Value = '4.4 °C'

# Row 34
$ws.Range("E34").Value = '2026-02-10 21:49:35'
$ws.Range("H34").Value = "'80%"
$ws.Range("O34").Value = '3.7 °C'

# Row 35
$ws.Range("E35").Value = '2026-02-10 21:49:37'

# Row 36
$ws.Range("E36").Value = '2026-02-10 21:49:40'
$ws.Range("H36").Value = "'94%"
$ws.Range("I36").Value = '5.6 mm'
$ws.Range("M36").Value = '13.0 °C 21:24 TU'
$ws.Range("O36").Value = '10.1 °C'

# Row 37
$ws.Range("E37").Value = '2026-02-10 21:49:42'
$ws.Range("J37").Value = '1005.5 hPa'

# Row 38
$ws.Range("E38").Value = '2026-02-10 21:49:45'
$ws.Range("K38").Value = '7.9 MJ/m2'
$ws.Range("O38").Value = '11.2 °C'

# Row 39
$ws.Range("E39").Value = '2026-02-10 21:49:47'
$ws.Range("H39").Value = "'81%"

# Row 40
$ws.Range("E40").Value = '2026-02-10 21:49:49'
$ws.Range("J40").Value = '1006.8 hPa'

# Row 41
$ws.Range("E41").Value = '2026-02-10 21:49:52'
$ws.Range("H41").Value = "'80%"
$ws.Range("J41").Value = '1004.4 hPa'
$ws.Range("O41").Value = '14.7 °C'

# Row 42
$ws.Range("E42").Value = '2026-02-10 21:49:54'
$ws.Range("I42").Value = '1.0 mm'

# Row 43
$ws.Range("E43").Value = '2026-02-10 21:49:57'
$ws.Range("K43").Value = '9.0 MJ/m2'
$ws.Range("O43").Value = '10.0 °C'

# Row 44
$ws.Range("E44").Value = '2026-02-10 21:49:59'
$ws.Range("I44").Value = '29.5 mm'
$ws.Range("O44").Value = '0.4 °C'

# Row 45
$ws.Range("E45").Value = '2026-02-10 21:50:01'
$ws.Range("J45").Value = '1005.6 hPa'

# Row 46
$ws.Range("E46").Value = '2026-02-10 21:50:04'
$ws.Range("H46").Value = "'79%"
$ws.Range("J46").Value = '1005.7 hPa'
$ws.Range("O46").Value = '14.9 °C'
